$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Il13"
$ws.Cells.Item(2, 3).Value = "Il13ra1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1174156666666667
$ws.Cells.Item(2, 8).Value = 0.352247
$ws.Cells.Item(2, 9).Value = 0.6540874079906115
$ws.Cells.Item(2, 10).Value = 0.7393359457808691
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 8.280552
$ws.Cells.Item(2, 14).Value = 16.561104
$ws.Cells.Item(2, 15).Value = 0.06375227685066993
$ws.Cells.Item(2, 16).Value = 0.04445819681279333
$ws.Cells.Item(2, 17).Value = 0.972266533448
$ws.Cells.Item(2, 18).Value = 5.833599200688
$ws.Cells.Item(2, 19).Value = 0.04169956151875456
$ws.Cells.Item(2, 20).Value = 0.03286954298829858

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Il13"
$ws.Cells.Item(3, 3).Value = "Il13ra1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1174156666666667
$ws.Cells.Item(3, 8).Value = 0.352247
$ws.Cells.Item(3, 9).Value = 0.6540874079906115
$ws.Cells.Item(3, 10).Value = 0.7393359457808691
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 27.78704333333333
$ws.Cells.Item(3, 14).Value = 83.36113
$ws.Cells.Item(3, 15).Value = 0.2139334768320069
$ws.Cells.Item(3, 16).Value = 0.2237825161943824
$ws.Cells.Item(3, 17).Value = 3.262634217678889
$ws.Cells.Item(3, 18).Value = 29.36370795911
$ws.Cells.Item(3, 19).Value = 0.139931193343467
$ws.Cells.Item(3, 20).Value = 0.1654504582597963

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Il13"
$ws.Cells.Item(4, 3).Value = "Il13ra1"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.1174156666666667
$ws.Cells.Item(4, 8).Value = 0.352247
$ws.Cells.Item(4, 9).Value = 0.6540874079906115
$ws.Cells.Item(4, 10).Value = 0.7393359457808691
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 26.34514033333334
$ws.Cells.Item(4, 14).Value = 79.03542100000001
$ws.Cells.Item(4, 15).Value = 0.2028322121762435
$ws.Cells.Item(4, 16).Value = 0.2121701730754169
$ws.Cells.Item(4, 17).Value = 3.093332215665222
$ws.Cells.Item(4, 18).Value = 27.839989940987
$ws.Cells.Item(4, 19).Value = 0.1326699959193609
$ws.Cells.Item(4, 20).Value = 0.1568650355772041

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Il13"
$ws.Cells.Item(5, 3).Value = "Il13ra1"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.1174156666666667
$ws.Cells.Item(5, 8).Value = 0.352247
$ws.Cells.Item(5, 9).Value = 0.6540874079906115
$ws.Cells.Item(5, 10).Value = 0.7393359457808691
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 8.8689915
$ws.Cells.Item(5, 14).Value = 17.737983
$ws.Cells.Item(5, 15).Value = 0.06828269437764999
$ws.Cells.Item(5, 16).Value = 0.04761752231469486
$ws.Cells.Item(5, 17).Value = 1.0413585496335
$ws.Cells.Item(5, 18).Value = 6.248151297801
$ws.Cells.Item(5, 19).Value = 0.04466285057609218
$ws.Cells.Item(5, 20).Value = 0.03520534589627657

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Il13"
$ws.Cells.Item(6, 3).Value = "Il13ra1"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.1174156666666667
$ws.Cells.Item(6, 8).Value = 0.352247
$ws.Cells.Item(6, 9).Value = 0.6540874079906115
$ws.Cells.Item(6, 10).Value = 0.7393359457808691
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 51.00428400000001
$ws.Cells.Item(6, 14).Value = 153.012852
$ws.Cells.Item(6, 15).Value = 0.3926838735072487
$ws.Cells.Item(6, 16).Value = 0.4107621985287224
$ws.Cells.Item(6, 17).Value = 5.988702008716
$ws.Cells.Item(6, 18).Value = 53.898318078444
$ws.Cells.Item(6, 19).Value = 0.2568495769820694
$ws.Cells.Item(6, 20).Value = 0.3036912585402621

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Il13"
$ws.Cells.Item(7, 3).Value = "Il13ra1"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.1174156666666667
$ws.Cells.Item(7, 8).Value = 0.352247
$ws.Cells.Item(7, 9).Value = 0.6540874079906115
$ws.Cells.Item(7, 10).Value = 0.7393359457808691
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.600361666666665
$ws.Cells.Item(7, 14).Value = 22.801085
$ws.Cells.Item(7, 15).Value = 0.05851546625618104
$ws.Cells.Item(7, 16).Value = 0.0612093930739901
$ws.Cells.Item(7, 17).Value = 0.8924015319994442
$ws.Cells.Item(7, 18).Value = 8.031613787994999
$ws.Cells.Item(7, 19).Value = 0.03827422965086755
$ws.Cells.Item(7, 20).Value = 0.04525430451903145

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Il13"
$ws.Cells.Item(8, 3).Value = "Il13ra1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.5
$ws.Cells.Item(8, 7).Value = 0.062095
$ws.Cells.Item(8, 8).Value = 0.12419
$ws.Cells.Item(8, 9).Value = 0.3459125920093885
$ws.Cells.Item(8, 10).Value = 0.2606640542191307
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 8.280552
$ws.Cells.Item(8, 14).Value = 16.561104
$ws.Cells.Item(8, 15).Value = 0.06375227685066993
$ws.Cells.Item(8, 16).Value = 0.04445819681279333
$ws.Cells.Item(8, 17).Value = 0.51418087644
$ws.Cells.Item(8, 18).Value = 2.05672350576
$ws.Cells.Item(8, 19).Value = 0.02205271533191537
$ws.Cells.Item(8, 20).Value = 0.01158865382449475

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Il13"
$ws.Cells.Item(9, 3).Value = "Il13ra1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.5
$ws.Cells.Item(9, 7).Value = 0.062095
$ws.Cells.Item(9, 8).Value = 0.12419
$ws.Cells.Item(9, 9).Value = 0.3459125920093885
$ws.Cells.Item(9, 10).Value = 0.2606640542191307
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 27.78704333333333
$ws.Cells.Item(9, 14).Value = 83.36113
$ws.Cells.Item(9, 15).Value = 0.2139334768320069
$ws.Cells.Item(9, 16).Value = 0.2237825161943824
$ws.Cells.Item(9, 17).Value = 1.725436455783333
$ws.Cells.Item(9, 18).Value = 10.3526187347
$ws.Cells.Item(9, 19).Value = 0.07400228348853999
$ws.Cells.Item(9, 20).Value = 0.05833205793458598

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Il13"
$ws.Cells.Item(10, 3).Value = "Il13ra1"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.5
$ws.Cells.Item(10, 7).Value = 0.062095
$ws.Cells.Item(10, 8).Value = 0.12419
$ws.Cells.Item(10, 9).Value = 0.3459125920093885
$ws.Cells.Item(10, 10).Value = 0.2606640542191307
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 26.34514033333334
$ws.Cells.Item(10, 14).Value = 79.03542100000001
$ws.Cells.Item(10, 15).Value = 0.2028322121762435
$ws.Cells.Item(10, 16).Value = 0.2121701730754169
$ws.Cells.Item(10, 17).Value = 1.635901488998333
$ws.Cells.Item(10, 18).Value = 9.815408933990001
$ws.Cells.Item(10, 19).Value = 0.07016221625688264
$ws.Cells.Item(10, 20).Value = 0.05530513749821284

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Il13"
$ws.Cells.Item(11, 3).Value = "Il13ra1"
$ws.Cells.Item(11, 4).Value = "MuSCs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.5
$ws.Cells.Item(11, 7).Value = 0.062095
$ws.Cells.Item(11, 8).Value = 0.12419
$ws.Cells.Item(11, 9).Value = 0.3459125920093885
$ws.Cells.Item(11, 10).Value = 0.2606640542191307
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 8.8689915
$ws.Cells.Item(11, 14).Value = 17.737983
$ws.Cells.Item(11, 15).Value = 0.06828269437764999
$ws.Cells.Item(11, 16).Value = 0.04761752231469486
$ws.Cells.Item(11, 17).Value = 0.5507200271925
$ws.Cells.Item(11, 18).Value = 2.20288010877
$ws.Cells.Item(11, 19).Value = 0.02361984380155781
$ws.Cells.Item(11, 20).Value = 0.01241217641841829

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Il13"
$ws.Cells.Item(12, 3).Value = "Il13ra1"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.5
$ws.Cells.Item(12, 7).Value = 0.062095
$ws.Cells.Item(12, 8).Value = 0.12419
$ws.Cells.Item(12, 9).Value = 0.3459125920093885
$ws.Cells.Item(12, 10).Value = 0.2606640542191307
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 51.00428400000001
$ws.Cells.Item(12, 14).Value = 153.012852
$ws.Cells.Item(12, 15).Value = 0.3926838735072487
$ws.Cells.Item(12, 16).Value = 0.4107621985287224
$ws.Cells.Item(12, 17).Value = 3.16711101498
$ws.Cells.Item(12, 18).Value = 19.00266608988
$ws.Cells.Item(12, 19).Value = 0.1358342965251793
$ws.Cells.Item(12, 20).Value = 0.1070709399884602

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Il13"
$ws.Cells.Item(13, 3).Value = "Il13ra1"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.5
$ws.Cells.Item(13, 7).Value = 0.062095
$ws.Cells.Item(13, 8).Value = 0.12419
$ws.Cells.Item(13, 9).Value = 0.3459125920093885
$ws.Cells.Item(13, 10).Value = 0.2606640542191307
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 7.600361666666665
$ws.Cells.Item(13, 14).Value = 22.801085
$ws.Cells.Item(13, 15).Value = 0.05851546625618104
$ws.Cells.Item(13, 16).Value = 0.0612093930739901
$ws.Cells.Item(13, 17).Value = 0.4719444576916665
$ws.Cells.Item(13, 18).Value = 2.831666746149999
$ws.Cells.Item(13, 19).Value = 0.02024123660531349
$ws.Cells.Item(13, 20).Value = 0.01595508855495864
